$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 79, shifting existing rows 79:133 down to 80:134.
$ws.Rows("79:79").Insert()

# Populate the new row 79 with the new data (constant columns copied from the
# surrounding rows, variable columns set to the new values from the diff).
$ws.Range("A79").Value = 3
$ws.Range("B79").Value = "Femacal de La Calera"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44596
$ws.Range("E79").Value = 5
$ws.Range("F79").Value = 100112052
$ws.Range("G79").Value = "Albahaca"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 165
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = 4242
$ws.Range("N79").Value = "$/docena de matas"
$ws.Range("O79").Value = "Provincia de Quillota"
$ws.Range("P79").Value = 707
$ws.Range("Q79").Value = 6
$ws.Range("R79").Value = "Hortaliza"

# Match the date-cell style used by the rest of column D.
$ws.Range("D79").NumberFormat = $ws.Range("D80").NumberFormat
